$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 26; $row++) {
    $ws.Range("B$row").Value = 0.9999549410293126
    $ws.Range("C$row").Value = 0.9990639202211722
    $ws.Range("D$row").Value = 0.999950352963198
    $ws.Range("E$row").Value = 0.9999999999995333
    $ws.Range("F$row").Value = 0.9999587938038491
    $ws.Range("G$row").Value = 0.00004206057127643839
    $ws.Range("H$row").Value = 0.0008737893843822983
    $ws.Range("I$row").Value = 0.00007622417229076227
    $ws.Range("J$row").Value = 0.0000000000001464511707968461
    $ws.Range("K$row").Value = 0.00003811208621860672
    $ws.Range("L$row").Value = 0.0004236862227404525
    $ws.Range("M$row").Value = 0.006485412190172526
    $ws.Range("N$row").Value = 1.0000514959665
    $ws.Range("O$row").Value = 0.00676150961376217
    $ws.Range("P$row").Value = 110.1527996107309
    $ws.Range("Q$row").Value = 165.0022117297999
}
